$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2022" column (Q) is appended after the existing "2021" column (P),
# for both the header row (4) and the data row (5). Copy the formatting
# from the adjoining column P so the new cells match the existing table
# styling, then fill in the new values.

$ws.Range("P4").Copy() | Out-Null
$ws.Range("Q4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("Q4").Value = 2022

$ws.Range("P5").Copy() | Out-Null
$ws.Range("Q5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("Q5").Value = 64.2

$excel.CutCopyMode = 0

# The workbook's last saved selection moves to R4 (just past the new
# column) once the edit is made.
$ws.Range("R4").Select() | Out-Null
